$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.0164
$ws.Range("G2").Value = -0.04008438818565401
$ws.Range("H2").Value = -0.04008438818565401
$ws.Range("I2").Value = -0.1054852320675106
$ws.Range("J2").Value = -0.1054852320675106
$ws.Range("K2").Value = -0.08
$ws.Range("L2").Value = -0.08438818565400845
$ws.Range("O2").Value = -0
$ws.Range("R2").Value = -0
$ws.Range("U2").Value = 0.173
$ws.Range("V2").Value = 0.1130718954248366
$ws.Range("W2").Value = -0.05633802816901409
$ws.Range("X2").Value = 0.06007044195043747
$ws.Range("Y2").Value = -0.1164084701194515
$ws.Range("Z2").Value = 0.6975717439293599
$ws.Range("AA2").Value = -0.07358351729212657
$ws.Range("AB2").Value = 0.05969669782470365
$ws.Range("AC2").Value = -0.1332802151168302
$ws.Range("AD2").Value = 0.021
$ws.Range("AF2").Value = 0.021
$ws.Range("AG2").Value = -0.152
$ws.Range("AH2").Value = 0.01353965183752418
$ws.Range("AI2").Value = 0.01389808074123097
$ws.Range("AJ2").Value = -0.1103047895500726
$ws.Range("AK2").Value = -0.1136023916292974
$ws.Range("AN2").Value = -0.2658227848101266
$ws.Range("AO2").Value = -100
$ws.Range("AP2").Value = 1.924050632911392
$ws.Range("AQ2").Value = -100

# Row 3
$ws.Range("D3").Value = -0.0164
$ws.Range("G3").Value = -0.04008438818565401
$ws.Range("H3").Value = -0.04008438818565401
$ws.Range("I3").Value = -0.1054852320675106
$ws.Range("J3").Value = -0.1054852320675106
$ws.Range("K3").Value = -0.08
$ws.Range("L3").Value = -0.08438818565400845
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 0.173
$ws.Range("V3").Value = 0.1130718954248366
$ws.Range("W3").Value = -0.05633802816901409
$ws.Range("X3").Value = 0.06007044195043747
$ws.Range("Y3").Value = -0.1164084701194515
$ws.Range("Z3").Value = 0.6975717439293599
$ws.Range("AA3").Value = -0.07358351729212657
$ws.Range("AB3").Value = 0.05969669782470365
$ws.Range("AC3").Value = -0.1332802151168302
$ws.Range("AD3").Value = 0.021
$ws.Range("AF3").Value = 0.021
$ws.Range("AG3").Value = -0.152
$ws.Range("AH3").Value = 0.01353965183752418
$ws.Range("AI3").Value = 0.01389808074123097
$ws.Range("AJ3").Value = -0.1103047895500726
$ws.Range("AK3").Value = -0.1136023916292974
$ws.Range("AN3").Value = -0.2658227848101266
$ws.Range("AO3").Value = -100
$ws.Range("AP3").Value = 1.924050632911392
$ws.Range("AQ3").Value = -100
